# Config_AR.xlsx — "Suppression du globalHandler. E2E gestion des nouvelles
# transactions. Reste le CR"
#
# Functional change: on the "CustomSettings" sheet, a new status-mapping
# row is inserted right above the existing "StatusRPA_OK_Gaps" row (the
# StatusRPA_* block that starts at row 31). The new row carries:
#   A = "StatusRPA_KO_Other"
#   B = "Autre erreur"
# Inserting the row pushes every row below it (incl. the trailing blank
# rows) down by one, which is exactly what Excel's Rows.Insert does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CustomSettings")
$ws.Activate()

# Insert a brand-new row above row 36 (the "StatusRPA_OK_Gaps" row),
# pushing it (and everything below) down to row 37.
$ws.Rows("36:36").Insert()

# The freshly inserted row 36 starts out with default formatting; pick up
# the same per-cell styles used by its neighbours (A:s84 / B:s73 / C:s65 /
# D:s63), matching the rest of the StatusRPA_* block above it, by copying
# formats from the row that used to be 36 and is now 37.
$ws.Range("A37").Copy()
$ws.Range("A36").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B37").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("C37").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("D37").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Rows("36:36").RowHeight = $ws.Rows("37:37").RowHeight
$excel.CutCopyMode = $false

$ws.Cells.Item(36, 1).Value = "StatusRPA_KO_Other"
$ws.Cells.Item(36, 2).Value = "Autre erreur"

# Restore the view state: selection moved to B34.
$ws.Range("B34").Select()
